{"js": "// ASNN-447: append 4 new TOC rows to the end of the (single) table,\n// mirroring the existing row formatting (Trebuchet MS, bold, blue\n// underlined title / right-aligned bold page number).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newRows = [\n  [\"SV: LATE SUBMISSION ASSIGNMENT LIST\", \"143\"],\n  [\"SV: VIEW DETAILS AND SUBMIT LATE SUBMISSION\", \"145\"],\n  [\"SV: LATE ASSIGNMENT LIST CONFIRMATION\", \"147\"],\n  [\"SV: VIEW SUBMISSION-LATE SUBMISSION\", \"149\"],\n];\n\ntable.addRows(\"End\", newRows.length, newRows);\nawait context.sync();\n", "ps1": "# ASNN-447: append 4 new TOC rows to the end of the (single) table,\n# mirroring the existing row formatting (Trebuchet MS, bold, blue\n# underlined title / right-aligned bold page number).\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newRows = @(\n    @(\"SV: LATE SUBMISSION ASSIGNMENT LIST\", \"143\"),\n    @(\"SV: VIEW DETAILS AND SUBMIT LATE SUBMISSION\", \"145\"),\n    @(\"SV: LATE ASSIGNMENT LIST CONFIRMATION\", \"147\"),\n    @(\"SV: VIEW SUBMISSION-LATE SUBMISSION\", \"149\")\n)\n\nforeach ($rowData in $newRows) {\n    $row = $t.Rows.Add()\n    $row.Cells.Item(1).Range.Text = $rowData[0]\n    $row.Cells.Item(2).Range.Text = $rowData[1]\n}\n"}
